$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 331 (pushes existing rows 331..351 down to 332..352)
$ws.Rows.Item(331).Insert()

# Populate the new row 331 with the weekly record.
# Columns A,B,C,E,F,G,H,I,R are constant for this market/product block.
$ws.Cells.Item(331, 1).Value = 10
$ws.Cells.Item(331, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(331, 3).Value = "La Araucanía"
$ws.Cells.Item(331, 4).Value = 44610
$ws.Cells.Item(331, 5).Value = 9
$ws.Cells.Item(331, 6).Value = 100114014
$ws.Cells.Item(331, 7).Value = "Betarraga"
$ws.Cells.Item(331, 8).Value = "Sin especificar"
$ws.Cells.Item(331, 9).Value = "Primera"
$ws.Cells.Item(331, 10).Value = 50
$ws.Cells.Item(331, 11).Value = 8000
$ws.Cells.Item(331, 12).Value = 8000
$ws.Cells.Item(331, 13).Value = 8000
$ws.Cells.Item(331, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(331, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(331, 16).Value = 667
$ws.Cells.Item(331, 17).Value = 12
$ws.Cells.Item(331, 18).Value = "Hortaliza"
